$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2310.818
$ws.Range("J19").Value = 2052.25
$ws.Range("L19").Value = 2052.25
$ws.Range("N19").Value = -2402.25
$ws.Range("H29").Value = 5816.6665
$ws.Range("I29").Value = 1300.3334
$ws.Range("K29").Value = 3901.0002
$ws.Range("M29").Value = -3620.0002
$ws.Range("H42").Value = 134.14285
$ws.Range("J42").Value = 90
$ws.Range("L42").Value = 270
$ws.Range("N42").Value = -730
$ws.Range("H118").Value = 318.57144
$ws.Range("I118").Value = 321.66666
$ws.Range("J118").Value = 300
$ws.Range("K118").Value = 964.9999799999999
$ws.Range("L118").Value = 900
$ws.Range("M118").Value = 692.0000200000001
$ws.Range("N118").Value = -4214
$ws.Range("H125").Value = 2238.25
$ws.Range("I125").Value = 2160.875
$ws.Range("K125").Value = 19447.875
$ws.Range("M125").Value = -16987.875
$ws.Range("H132").Value = 2577.3076
$ws.Range("I132").Value = 2250.4167
$ws.Range("K132").Value = 6751.250100000001
$ws.Range("M132").Value = -4221.250100000001
$ws.Range("H135").Value = 14199.889
$ws.Range("I135").Value = 2633.3333
$ws.Range("J135").Value = 37333
$ws.Range("K135").Value = 23699.9997
$ws.Range("L135").Value = 335997
$ws.Range("M135").Value = -21164.9997
$ws.Range("N135").Value = -341067
$ws.Range("H138").Value = 2941.443
$ws.Range("J138").Value = 3264.2646
$ws.Range("L138").Value = 9792.793799999999
$ws.Range("N138").Value = -20072.7938
$ws.Range("H139").Value = 88000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 88000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 88000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -98280
$ws.Range("H141").Value = 5039.364
$ws.Range("I141").Value = 5039.364
$ws.Range("K141").Value = 15118.092
$ws.Range("M141").Value = -9938.091999999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31263202
$ws.Range("I32").Value = 31263202
$ws.Range("K32").Value = 31263202
$ws.Range("M32").Value = -31262915
$ws.Range("H45").Value = 2500
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H61").Value = 23863982
$ws.Range("I61").Value = 45459024
$ws.Range("J61").Value = 109436
$ws.Range("K61").Value = 45459024
$ws.Range("L61").Value = 109436
$ws.Range("M61").Value = -45458812
$ws.Range("N61").Value = -109860
$ws.Range("H63").Value = 3077.8
$ws.Range("I63").Value = 3077.8
$ws.Range("K63").Value = 3077.8
$ws.Range("M63").Value = -2391.8
$ws.Range("H66").Value = 3077.8
$ws.Range("I66").Value = 3077.8
$ws.Range("K66").Value = 15389
$ws.Range("M66").Value = -11957
$ws.Range("H74").Value = 13898785
$ws.Range("I74").Value = 19231828
$ws.Range("K74").Value = 19231828
$ws.Range("M74").Value = -19230954
$ws.Range("H77").Value = 13898785
$ws.Range("I77").Value = 19231828
$ws.Range("K77").Value = 96159140
$ws.Range("M77").Value = -96154772
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2000
$ws.Range("K122").Value = 6000
$ws.Range("M122").Value = -3550
$ws.Range("H132").Value = 6607.864
$ws.Range("I132").Value = 3676.125
$ws.Range("K132").Value = 11028.375
$ws.Range("M132").Value = -8498.375
$ws.Range("H136").Value = 23863982
$ws.Range("I136").Value = 45459024
$ws.Range("J136").Value = 109436
$ws.Range("K136").Value = 136377072
$ws.Range("L136").Value = 328308
$ws.Range("M136").Value = -136374522
$ws.Range("N136").Value = -333408
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4535.3335
$ws.Range("I86").Value = 4191.143
$ws.Range("J86").Value = 5740
$ws.Range("K86").Value = 4191.143
$ws.Range("L86").Value = 5740
$ws.Range("M86").Value = -3068.143
$ws.Range("N86").Value = -7986
$ws.Range("H89").Value = 4535.3335
$ws.Range("I89").Value = 4191.143
$ws.Range("J89").Value = 5740
$ws.Range("K89").Value = 20955.715
$ws.Range("L89").Value = 28700
$ws.Range("M89").Value = -15339.715
$ws.Range("N89").Value = -39932
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 835118.4
$ws.Range("I31").Value = 1791.5
$ws.Range("J31").Value = 1946220.9
$ws.Range("K31").Value = 1791.5
$ws.Range("L31").Value = 1946220.9
$ws.Range("M31").Value = -1496.5
$ws.Range("N31").Value = -1946810.9
$ws.Range("H34").Value = 835118.4
$ws.Range("I34").Value = 1791.5
$ws.Range("J34").Value = 1946220.9
$ws.Range("K34").Value = 1791.5
$ws.Range("L34").Value = 1946220.9
$ws.Range("M34").Value = -1589.5
$ws.Range("N34").Value = -1946624.9
$ws.Range("H124").Value = 36698.668
$ws.Range("J124").Value = 36698.668
$ws.Range("L124").Value = 36698.668
$ws.Range("N124").Value = -41608.668
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 629.8
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168
$ws.Range("H39").Value = 338719
$ws.Range("J39").Value = 376499.75
$ws.Range("L39").Value = 1129499.25
$ws.Range("N39").Value = -1130087.25
$ws.Range("H55").Value = 4000
$ws.Range("J55").Value = 4750
$ws.Range("L55").Value = 14250
$ws.Range("N55").Value = -14604
$ws.Range("H113").Value = 1299.6666
$ws.Range("I113").Value = 494.66666
$ws.Range("K113").Value = 1483.99998
$ws.Range("M113").Value = 686.0000199999999
$ws.Range("H123").Value = 5666.1665
$ws.Range("I123").Value = 4999.5
$ws.Range("K123").Value = 14998.5
$ws.Range("M123").Value = -12548.5
$ws.Range("H131").Value = 9819
$ws.Range("I131").Value = 9342
$ws.Range("J131").Value = 11250
$ws.Range("K131").Value = 28026
$ws.Range("L131").Value = 33750
$ws.Range("M131").Value = -22986
$ws.Range("N131").Value = -43830
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 62502948
$ws.Range("I132").Value = 71431510
$ws.Range("K132").Value = 214294530
$ws.Range("M132").Value = -214292000
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2817
$ws.Range("J16").Value = 2000
$ws.Range("L16").Value = 2000
$ws.Range("N16").Value = -2340
$ws.Range("H22").Value = 4017
$ws.Range("J22").Value = 2650
$ws.Range("L22").Value = 2650
$ws.Range("N22").Value = -3240
$ws.Range("H27").Value = 4017
$ws.Range("J27").Value = 2650
$ws.Range("L27").Value = 2650
$ws.Range("N27").Value = -2864
$ws.Range("H68").Value = 3421.75
$ws.Range("I68").Value = 3161.3333
$ws.Range("J68").Value = 3812.375
$ws.Range("K68").Value = 3161.3333
$ws.Range("L68").Value = 3812.375
$ws.Range("M68").Value = -2412.3333
$ws.Range("N68").Value = -5310.375
$ws.Range("H71").Value = 3421.75
$ws.Range("I71").Value = 3161.3333
$ws.Range("J71").Value = 3812.375
$ws.Range("K71").Value = 15806.6665
$ws.Range("L71").Value = 19061.875
$ws.Range("M71").Value = -12062.6665
$ws.Range("N71").Value = -26549.875
$ws.Range("H122").Value = 5899.1055
$ws.Range("I122").Value = 5206.9165
$ws.Range("K122").Value = 15620.7495
$ws.Range("M122").Value = -13170.7495
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3596.1316
$ws.Range("I122").Value = 2419.5833
$ws.Range("K122").Value = 7258.749899999999
$ws.Range("M122").Value = -4808.749899999999
$ws.Range("H138").Value = 191499.5
$ws.Range("J138").Value = 191499.5
$ws.Range("L138").Value = 191499.5
$ws.Range("N138").Value = -201779.5
